{"js": "// Replace the 15 lattice-multiplication exercise cells with a new set of\n// problems. Each table cell holds a single run made of five <w:t> pieces\n// joined by <w:br/> line breaks:\n//   \"A x B\"\n//   \"  <tens of B>    <ones of B>\"\n//   \"  ----\"\n//   \"<tens of A>|    |\"\n//   \"<ones of A>|    |\"\n// (vertical tab, \\u000b, is how Office.js exposes a <w:br/> inside text.)\n\nconst VT = \"\\u000b\";\n\nfunction buildCellText(a, b) {\n  const aTens = Math.floor(a / 10).toString();\n  const aOnes = (a % 10).toString();\n  const bTens = Math.floor(b / 10).toString();\n  const bOnes = (b % 10).toString();\n  return (\n    a + \" x \" + b + VT +\n    \"  \" + bTens + \"    \" + bOnes + VT +\n    \"  ----\" + VT +\n    aTens + \"|    |\" + VT +\n    aOnes + \"|    |\"\n  );\n}\n\n// New problems, in row-major order (5 rows x 3 columns), taken from the\n// target document.\nconst problems = [\n  [91, 91], [61, 91], [14, 24],\n  [37, 50], [92, 21], [68, 45],\n  [79, 83], [75, 62], [18, 79],\n  [80, 22], [71, 36], [59, 96],\n  [16, 60], [26, 45], [49, 12],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 3;\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const idx = r * colCount + c;\n    if (idx >= problems.length) continue;\n    const [a, b] = problems[idx];\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(buildCellText(a, b), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication exercise cells with a new set of\n# problems. Each table cell holds a single run made of five text pieces\n# joined by line-break characters (Chr(11), the same character a <w:br/>\n# round-trips to in the Range.Text string):\n#   \"A x B\"\n#   \"  <tens of B>    <ones of B>\"\n#   \"  ----\"\n#   \"<tens of A>|    |\"\n#   \"<ones of A>|    |\"\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$vt = [char]11\n\n# New problems, in row-major order (5 rows x 3 columns), taken from the\n# target document.\n$problems = @(\n    @(91, 91), @(61, 91), @(14, 24),\n    @(37, 50), @(92, 21), @(68, 45),\n    @(79, 83), @(75, 62), @(18, 79),\n    @(80, 22), @(71, 36), @(59, 96),\n    @(16, 60), @(26, 45), @(49, 12)\n)\n\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $problems.Count) { continue }\n        $pair = $problems[$idx]\n        $a = $pair[0]\n        $b = $pair[1]\n\n        $aTens = [int]([math]::Floor($a / 10))\n        $aOnes = $a % 10\n        $bTens = [int]([math]::Floor($b / 10))\n        $bOnes = $b % 10\n\n        $text = \"$a x $b\" + $vt + \"  $bTens    $bOnes\" + $vt + \"  ----\" + $vt + \"$aTens|    |\" + $vt + \"$aOnes|    |\"\n\n        $cell = $table.Cell($r, $c)\n        $rng = $cell.Range\n        $rng.End = $rng.End - 1\n        $rng.Text = $text\n\n        $idx = $idx + 1\n    }\n}\n"}
